# Remove the "high-res lighting" option from the localization sheet.
#
# Row 30 (A30="high-res lighting", B30="High-resolution color") is the
# row being removed. Deleting it as a whole row shifts every row below it
# up by one — which explains why the diff shows no visible change to the
# sheet content until row 55 onward (rows 55/56, 62/63, 71/72 and 76/77
# simply carry their own per-row height/style along as they move, which
# is what produces the apparent attribute "swap" in the diff), and why
# the final row (85) disappears, dimension shrinks by one row, and two
# shared strings are pruned.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A30").EntireRow.Delete()

# The conditional-formatting range that watches column C ("C2:<last row>")
# was anchored to the sheet's bottom row; since a row disappeared from
# inside it, its extent shrinks by one row too (C2:C1048576 -> C2:C1048575).
$cf = $ws.Range("C2:C1048576").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("C2:C1048575"))

# Leave the selection where the deleted row used to be, matching the
# author's resulting cursor position after the delete.
$ws.Range("A30").Select()
